# PHQ parameter calculation done
# Sheet2 (3rd tab, sheet3.xml) gets two new "0D Tmax"/"0D Tavg" columns inserted
# after "exp Tavg" (duplicating the exp Tmax/Tavg values), plus two new
# "fitted Tmax"/"fitted Tavg" columns appended at the end with freshly
# computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- insert two new columns before the old "kh" column (E) -----------------
$ws.Range("E1:F1").EntireColumn.Insert()

# --- new header E1/F1: "0D Tmax" / "0D Tavg" --------------------------------
$ws.Range("E1").Value = "0D Tmax"
$ws.Range("F1").Value = "0D Tavg"

# --- new data for E2:F17 (same values as exp Tmax/exp Tavg, columns C/D) ---
$expTmax = @(954.4, 1140.3, 1324.3, 1480.1, 1627.5, 1827.8, 1947.8, 2025.6, 1018.5, 1235.8, 1453.4, 1748.5, 1988.1, 1140.4, 1487.3, 1816.2)
$expTavg = @(767.4, 792.7, 820.8, 876, 908.8, 982.6, 998.4, 1003.2, 769.6, 821.1, 917.3, 979.1, 1014, 800.4, 899.3, 990.1)

for ($i = 0; $i -lt $expTmax.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $expTmax[$i]
    $ws.Cells.Item($r, 6).Value = $expTavg[$i]
}

# --- new header K1/L1: "fitted Tmax" / "fitted Tavg" ------------------------
$ws.Range("K1").Value = "fitted Tmax"
$ws.Range("L1").Value = "fitted Tavg"

# --- new data for K2:L17 (freshly fitted model output) ----------------------
$fitTmax = @(951.6, 1111.6, 1271.3, 1423.6, 1580.6, 1747.9, 1927.7, 2121, 1018.7, 1253.9, 1543.1, 1809.5, 2081, 1140.3, 1476.5, 1802.2)
$fitTavg = @(765.7, 797.7, 833.3, 870.4, 909.2, 950.1, 993.5, 1039.6, 769.8, 825.8, 891.4, 958.9, 1028.7, 800.3, 899.4, 993.3)

for ($i = 0; $i -lt $fitTmax.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 11).Value = $fitTmax[$i]
    $ws.Cells.Item($r, 12).Value = $fitTavg[$i]
}

# --- column widths for the new fitted-value columns -------------------------
$ws.Columns("K").ColumnWidth = 11.25
$ws.Columns("L").ColumnWidth = 10.625

# --- selection moves to O9 ---------------------------------------------------
$ws.Range("O9").Select() | Out-Null

# --- page setup: paper size / orientation -----------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
